$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.877.64"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "3.829.66"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "704.30"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.46"
$ws.Range("E6").Value = "  -1.49%  "

$ws.Range("D7").Value = "3.829.08"
$ws.Range("E7").Value = "  -0.65%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.43"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.51"
$ws.Range("E14").Value = "  -0.42%  "

$ws.Range("D15").Value = "4.475.53"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "3.885.40"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "70.992.92"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "494.60"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.32"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000144"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.60"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.08"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.09"
$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.41"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.31"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("E34").Value = "  -3.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.17"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("D36").Value = "3.792.28"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.05"
$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.97"
$ws.Range("E41").Value = "  -1.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.30"
$ws.Range("E42").Value = "  -3.51%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000311"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.51"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "429.75"
$ws.Range("E47").Value = "  +4.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.85"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.72"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("E50").Value = "  -1.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.295"
$ws.Range("E51").Value = "  -2.81%  "
